$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-13 18:18:40"
$ws.Range("G2").Value = "109 cm"
$ws.Range("I2").Value = "1.8 mm"
$ws.Range("E3").Value = "2026-02-13 18:18:43"
$ws.Range("I3").Value = "5.5 mm"
$ws.Range("E4").Value = "2026-02-13 18:18:45"
$ws.Range("I4").Value = "5.1 mm"
$ws.Range("J4").Value = "995.2 hPa"
$ws.Range("E5").Value = "2026-02-13 18:18:48"
$ws.Range("G5").Value = "105 cm"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "81%"
$ws.Range("I5").Value = "0.5 mm"
$ws.Range("E6").Value = "2026-02-13 18:18:51"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "74%"
$ws.Range("I6").Value = "4.5 mm"
$ws.Range("J6").Value = "995.3 hPa"
$ws.Range("E7").Value = "2026-02-13 18:18:54"
$ws.Range("I7").Value = "16.7 mm"
$ws.Range("J7").Value = "995.5 hPa"
$ws.Range("L7").Value = "32.4 km/h - 289º 17:51 TU"
$ws.Range("O7").Value = "13.0 °C"
$ws.Range("E8").Value = "2026-02-13 18:18:56"
$ws.Range("I8").Value = "19.3 mm"
$ws.Range("J8").Value = "995.4 hPa"
$ws.Range("N8").Value = "7.2 °C 17:59 TU"
$ws.Range("O8").Value = "9.5 °C"
$ws.Range("E9").Value = "2026-02-13 18:18:59"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "75%"
$ws.Range("I9").Value = "3.4 mm"
$ws.Range("E10").Value = "2026-02-13 18:19:02"
$ws.Range("I10").Value = "17.1 mm"
$ws.Range("O10").Value = "9.0 °C"
$ws.Range("E11").Value = "2026-02-13 18:19:04"
$ws.Range("E12").Value = "2026-02-13 18:19:07"
$ws.Range("I12").Value = "4.6 mm"
$ws.Range("E13").Value = "2026-02-13 18:19:09"
$ws.Range("I13").Value = "5.9 mm"
$ws.Range("J13").Value = "998.2 hPa"
$ws.Range("K13").Value = "2.5 MJ/m2"
$ws.Range("O13").Value = "0.5 °C"
$ws.Range("E14").Value = "2026-02-13 18:19:12"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "83%"
$ws.Range("I14").Value = "18.9 mm"
$ws.Range("E15").Value = "2026-02-13 18:19:14"
$ws.Range("I15").Value = "2.7 mm"
$ws.Range("O15").Value = "9.5 °C"
$ws.Range("E16").Value = "2026-02-13 18:19:17"
$ws.Range("I16").Value = "12.6 mm"
$ws.Range("E17").Value = "2026-02-13 18:19:20"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "89%"
$ws.Range("I17").Value = "4.7 mm"
$ws.Range("E18").Value = "2026-02-13 18:19:22"
$ws.Range("I18").Value = "8.1 mm"
$ws.Range("J18").Value = "995.4 hPa"
$ws.Range("E19").Value = "2026-02-13 18:19:25"
$ws.Range("I19").Value = "12.5 mm"
$ws.Range("E20").Value = "2026-02-13 18:19:28"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "93%"
$ws.Range("I20").Value = "20.0 mm"
$ws.Range("E21").Value = "2026-02-13 18:19:30"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "91%"
$ws.Range("J21").Value = "998.2 hPa"
$ws.Range("E22").Value = "2026-02-13 18:19:33"
$ws.Range("E23").Value = "2026-02-13 18:19:36"
$ws.Range("G23").Value = "185 cm"
$ws.Range("I23").Value = "8.7 mm"
$ws.Range("E24").Value = "2026-02-13 18:19:39"
$ws.Range("J24").Value = "996.0 hPa"
$ws.Range("E25").Value = "2026-02-13 18:19:42"
$ws.Range("I25").Value = "8.7 mm"
$ws.Range("E26").Value = "2026-02-13 18:19:44"
$ws.Range("E27").Value = "2026-02-13 18:19:47"
$ws.Range("E28").Value = "2026-02-13 18:19:49"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "78%"
$ws.Range("I28").Value = "6.5 mm"
$ws.Range("J28").Value = "995.7 hPa"
$ws.Range("E29").Value = "2026-02-13 18:19:52"
$ws.Range("I29").Value = "12.9 mm"
$ws.Range("E30").Value = "2026-02-13 18:19:55"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "77%"
$ws.Range("I30").Value = "3.5 mm"
$ws.Range("J30").Value = "995.2 hPa"
$ws.Range("E31").Value = "2026-02-13 18:19:58"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "71%"
$ws.Range("I31").Value = "2.3 mm"
$ws.Range("J31").Value = "994.2 hPa"
$ws.Range("N31").Value = "8.1 °C 17:50 TU"
$ws.Range("E32").Value = "2026-02-13 18:20:00"
$ws.Range("I32").Value = "23.7 mm"
$ws.Range("L32").Value = "43.2 km/h - 297º 17:47 TU"
$ws.Range("E33").Value = "2026-02-13 18:20:02"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "89%"
$ws.Range("I33").Value = "4.8 mm"
$ws.Range("J33").Value = "997.2 hPa"
$ws.Range("E34").Value = "2026-02-13 18:20:05"
$ws.Range("G34").Value = "106 cm"
$ws.Range("I34").Value = "10.3 mm"
$ws.Range("E35").Value = "2026-02-13 18:20:08"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "74%"
$ws.Range("I35").Value = "7.5 mm"
$ws.Range("J35").Value = "995.9 hPa"
$ws.Range("L35").Value = "71.6 km/h - 255º 17:44 TU"
$ws.Range("N35").Value = "3.7 °C 17:55 TU"
$ws.Range("O35").Value = "6.2 °C"
$ws.Range("E36").Value = "2026-02-13 18:20:11"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "77%"
$ws.Range("I36").Value = "7.9 mm"
$ws.Range("J36").Value = "995.4 hPa"
$ws.Range("O36").Value = "10.7 °C"
$ws.Range("E37").Value = "2026-02-13 18:19:14"
$ws.Range("I37").Value = "11.9 mm"
$ws.Range("J37").Value = "997.2 hPa"
$ws.Range("E38").Value = "2026-02-13 18:20:16"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "77%"
$ws.Range("I38").Value = "12.2 mm"
$ws.Range("E39").Value = "2026-02-13 18:20:19"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "78%"
$ws.Range("I39").Value = "17.9 mm"
$ws.Range("E40").Value = "2026-02-13 18:20:22"
$ws.Range("G40").Value = "2 cm"
$ws.Range("J40").Value = "998.7 hPa"
$ws.Range("E41").Value = "2026-02-13 18:20:24"
$ws.Range("J41").Value = "995.5 hPa"
$ws.Range("L41").Value = "49.3 km/h - 279º 17:47 TU"
$ws.Range("E42").Value = "2026-02-13 18:20:27"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "84%"
$ws.Range("I42").Value = "9.1 mm"
$ws.Range("O42").Value = "11.1 °C"
$ws.Range("E43").Value = "2026-02-13 18:20:29"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "87%"
$ws.Range("I43").Value = "12.3 mm"
$ws.Range("E44").Value = "2026-02-13 18:20:32"
$ws.Range("I44").Value = "5.4 mm"
$ws.Range("E45").Value = "2026-02-13 18:20:35"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "62%"
$ws.Range("J45").Value = "993.7 hPa"
$ws.Range("E46").Value = "2026-02-13 18:20:38"
$ws.Range("J46").Value = "996.1 hPa"
